$d = $word.ActiveDocument

# 1. Fix wording: "3_ at the start of a layer name." -> "Layers starting with 3_."
$d.Content.Find.Execute("3_ at the start of a layer name.", $true, $false, $false, $false, $false, $true, 1, $false, "Layers starting with 3_.", 2) | Out-Null

# 2. Tidy spacing in the Examples list (remove stray spaces before commas).
$d.Content.Find.Execute("Examples: 3_1_2_4, 2_3_4_1, 1_1 , 2_3, 2_4, 1_3_1_2 , 1_3_5_2.", $true, $false, $false, $false, $false, $true, 1, $false, "Examples: 3_1_2_4, 2_3_4_1, 1_1, 2_3, 2_4, 1_3_1_2, 1_3_5_2.", 2) | Out-Null

# 3. Append the new "Deleting an Attached Strand" section at the end of the document.
$newParagraphs = @(
  @{Text='Deleting an Attached Strand'; Style='Heading1'},
  @{Text='Example: Deleting x_y (where y ≠ 1)'; Style='Heading2'},
  @{Text='Objective'; Style='Heading3'},
  @{Text='When deleting an attached strand x_y (where y is not equal to 1), the goal is to remove:'; Style='Normal'},
  @{Text='The Specific Attached Strand:'; Style='Heading4'},
  @{Text='Delete the attached strand x_y.'; Style='Normal'},
  @{Text='Associated Masks:'; Style='Heading4'},
  @{Text='Delete all mask layers that include x_y as a component.'; Style='Normal'},
  @{Text='Examples:'; Style='Normal'},
  @{Text='z_w_x_y: A mask layer where x_y is part of the sequence.'; Style='Normal'},
  @{Text='x_y_z_w: Another mask pattern including x_y.'; Style='Normal'},
  @{Text='Naming Patterns'; Style='Heading3'},
  @{Text='Direct Relationship:'; Style='Heading4'},
  @{Text='x_y itself is an attached strand to be deleted.'; Style='Normal'},
  @{Text='Complex Naming Patterns:'; Style='Heading4'},
  @{Text='z_w_x_y: Includes x_y within the layer name.'; Style='Normal'},
  @{Text='x_y_z_w: Another pattern with x_y as a component.'; Style='Normal'},
  @{Text='Excluded Layers:'; Style='Heading4'},
  @{Text='Layers that do not include x_y in these significant positions should remain unchanged.'; Style='Normal'}
)

foreach ($item in $newParagraphs) {
    $count = $d.Paragraphs.Count
    $last = $d.Paragraphs.Item($count)
    $last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.set_Style($item.Style)
    $newPara.Range.Text = $item.Text
}

